# Scheduled runner refresh: update cached market-price / profit figures
# (columns H..N = currentAveragePrice*, LevePrice*, LeveProfit* etc.)
# across the leve-profit sheets. Values below were pulled from the
# latest market snapshot; row identity (leve name/item) is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5415.364
$ws.Range("I32").Value = 3767.6
$ws.Range("J32").Value = 6788.5
$ws.Range("K32").Value = 3767.6
$ws.Range("L32").Value = 6788.5
$ws.Range("M32").Value = -3441.6
$ws.Range("N32").Value = -7440.5

$ws.Range("H53").Value = 2495.7693
$ws.Range("I53").Value = 419.25
$ws.Range("J53").Value = 3418.6667
$ws.Range("K53").Value = 419.25
$ws.Range("L53").Value = 3418.6667
$ws.Range("M53").Value = 217.75
$ws.Range("N53").Value = -4692.6667

$ws.Range("H70").Value = 56074.55
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21540

$ws.Range("H73").Value = 56074.55
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -22872

$ws.Range("H95").Value = 71882.664
$ws.Range("J95").Value = 71882.664
$ws.Range("L95").Value = 71882.664
$ws.Range("N95").Value = -77374.664

$ws.Range("H100").Value = 4095.875
$ws.Range("I100").Value = 1748.8
$ws.Range("J100").Value = 8007.6665
$ws.Range("K100").Value = 1748.8
$ws.Range("L100").Value = 8007.6665
$ws.Range("M100").Value = -1207.8
$ws.Range("N100").Value = -9089.666499999999

$ws.Range("H112").Value = 1285.4849
$ws.Range("I112").Value = 1095.2
$ws.Range("J112").Value = 1319.4642
$ws.Range("K112").Value = 3285.6
$ws.Range("L112").Value = 3958.3926
$ws.Range("M112").Value = -2177.6
$ws.Range("N112").Value = -6174.392599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3112.9412
$ws.Range("I32").Value = 2260.7969
$ws.Range("K32").Value = 2260.7969
$ws.Range("M32").Value = -1973.7969

$ws.Range("H61").Value = 4204.6763
$ws.Range("I61").Value = 3025.7585
$ws.Range("K61").Value = 3025.7585
$ws.Range("M61").Value = -2813.7585

$ws.Range("H74").Value = 23813518
$ws.Range("I74").Value = 25643864
$ws.Range("K74").Value = 25643864
$ws.Range("M74").Value = -25642990

$ws.Range("H77").Value = 23813518
$ws.Range("I77").Value = 25643864
$ws.Range("K77").Value = 128219320
$ws.Range("M77").Value = -128214952

$ws.Range("H102").Value = 2510
$ws.Range("I102").Value = 2250
$ws.Range("K102").Value = 2250
$ws.Range("M102").Value = -628

$ws.Range("H122").Value = 3741.7727
$ws.Range("I122").Value = 3388.182
$ws.Range("K122").Value = 10164.546
$ws.Range("M122").Value = -7714.545999999998

$ws.Range("H136").Value = 4204.6763
$ws.Range("I136").Value = 3025.7585
$ws.Range("K136").Value = 9077.2755
$ws.Range("M136").Value = -6527.2755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2880.1428
$ws.Range("I99").Value = 2686.889
$ws.Range("K99").Value = 2686.889
$ws.Range("M99").Value = -1188.889

$ws.Range("H105").Value = 13942.174
$ws.Range("I105").Value = 16863.924
$ws.Range("J105").Value = 10143.9
$ws.Range("K105").Value = 16863.924
$ws.Range("L105").Value = 10143.9
$ws.Range("M105").Value = -15116.924
$ws.Range("N105").Value = -13637.9

$ws.Range("H107").Value = 3437.5557
$ws.Range("I107").Value = 3539.6667
$ws.Range("K107").Value = 3539.6667
$ws.Range("M107").Value = -1619.6667

$ws.Range("H138").Value = 64801.11
$ws.Range("J138").Value = 64801.11
$ws.Range("L138").Value = 64801.11
$ws.Range("N138").Value = -75081.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3053
$ws.Range("I22").Value = 1387
$ws.Range("J22").Value = 3719.4
$ws.Range("K22").Value = 1387
$ws.Range("L22").Value = 3719.4
$ws.Range("M22").Value = -1037
$ws.Range("N22").Value = -4419.4

$ws.Range("H31").Value = 101630.82
$ws.Range("I31").Value = 9337
$ws.Range("K31").Value = 9337
$ws.Range("M31").Value = -9042

$ws.Range("H34").Value = 101630.82
$ws.Range("I34").Value = 9337
$ws.Range("K34").Value = 9337
$ws.Range("M34").Value = -9135

$ws.Range("H58").Value = 10126.272
$ws.Range("I58").Value = 3133.6667
$ws.Range("J58").Value = 12748.5
$ws.Range("K58").Value = 3133.6667
$ws.Range("L58").Value = 12748.5
$ws.Range("M58").Value = -2930.6667
$ws.Range("N58").Value = -13154.5

$ws.Range("H99").Value = 3928.4285
$ws.Range("I99").Value = 3666.3333
$ws.Range("K99").Value = 3666.3333
$ws.Range("M99").Value = -2168.3333

$ws.Range("H107").Value = 1515.35
$ws.Range("I107").Value = 1344.5625
$ws.Range("J107").Value = 2198.5
$ws.Range("K107").Value = 1344.5625
$ws.Range("L107").Value = 2198.5
$ws.Range("M107").Value = 575.4375
$ws.Range("N107").Value = -6038.5

$ws.Range("H126").Value = 3928.4285
$ws.Range("I126").Value = 3666.3333
$ws.Range("K126").Value = 10998.9999
$ws.Range("M126").Value = -8528.999899999999

$ws.Range("H136").Value = 10126.272
$ws.Range("I136").Value = 3133.6667
$ws.Range("J136").Value = 12748.5
$ws.Range("K136").Value = 9401.000100000001
$ws.Range("L136").Value = 38245.5
$ws.Range("M136").Value = -6851.000100000001
$ws.Range("N136").Value = -43345.5

$ws.Range("H141").Value = 274431.88
$ws.Range("J141").Value = 274431.88
$ws.Range("L141").Value = 274431.88
$ws.Range("N141").Value = -284791.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 154.375
$ws.Range("J2").Value = 169.83333
$ws.Range("L2").Value = 1018.99998
$ws.Range("N2").Value = -1244.99998

$ws.Range("H37").Value = 208647.47
$ws.Range("J37").Value = 208647.47
$ws.Range("L37").Value = 625942.41
$ws.Range("N37").Value = -626166.41

$ws.Range("H64").Value = 71434370
$ws.Range("I64").Value = 1000000000
$ws.Range("K64").Value = 3000000000
$ws.Range("M64").Value = -2999999730

$ws.Range("H67").Value = 71434370
$ws.Range("I67").Value = 1000000000
$ws.Range("K67").Value = 3000000000
$ws.Range("M67").Value = -2999999064

$ws.Range("H139").Value = 3543.625
$ws.Range("I139").Value = 2106
$ws.Range("J139").Value = 6706.4
$ws.Range("K139").Value = 6318
$ws.Range("L139").Value = 20119.2
$ws.Range("M139").Value = -1178
$ws.Range("N139").Value = -30399.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 49998
$ws.Range("J82").Value = 49998
$ws.Range("L82").Value = 49998
$ws.Range("N82").Value = -50764

$ws.Range("H85").Value = 49998
$ws.Range("J85").Value = 49998
$ws.Range("L85").Value = 49998
$ws.Range("N85").Value = -52650

$ws.Range("H107").Value = 1444.7727
$ws.Range("I107").Value = 1600.8948
$ws.Range("K107").Value = 4802.6844
$ws.Range("M107").Value = -2882.6844
